$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from existing "sum" header (G1) into new H1 header cell
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add data values for the new "Save" column
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0

$excel.CutCopyMode = $false

